$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 330 (shifts old rows 330..437 down to 331..438)
$ws.Rows.Item(330).Insert()

# Populate the newly inserted row 330 with the new weekly price record
$ws.Cells.Item(330, 1).Value = 6
$ws.Cells.Item(330, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(330, 3).Value = "Metropolitana"
$ws.Cells.Item(330, 4).Value = 44588
$ws.Cells.Item(330, 5).Value = 13
$ws.Cells.Item(330, 6).Value = 100112044
$ws.Cells.Item(330, 7).Value = "Perejil"
$ws.Cells.Item(330, 8).Value = "Sin especificar"
$ws.Cells.Item(330, 9).Value = "Primera"
$ws.Cells.Item(330, 10).Value = 240
$ws.Cells.Item(330, 11).Value = 14000
$ws.Cells.Item(330, 12).Value = 15000
$ws.Cells.Item(330, 13).Value = 14625
$ws.Cells.Item(330, 14).Value = "$/docena de atados"
$ws.Cells.Item(330, 15).Value = "Región Metropolitana"
$ws.Cells.Item(330, 16).Value = 4875
$ws.Cells.Item(330, 17).Value = 3
$ws.Cells.Item(330, 18).Value = "Hortaliza"
